# Regenerate merged AHB files
#
# 1. Rename the "_old" / "_new" header-name suffixes to "_FV2310" / "_FV2404"
# 2. Turn the used range into an Excel Table ("Table1")
# 3. Freeze the header row (row 1)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename header row ---------------------------------------------
$headers = @(
    "Segmentname_FV2310",
    "Segmentgruppe_FV2310",
    "Segment_FV2310",
    "Datenelement_FV2310",
    "Segment ID_FV2310",
    "Code_FV2310",
    "Qualifier_FV2310",
    "Beschreibung_FV2310",
    "Bedingungsausdruck_FV2310",
    "Bedingung_FV2310",
    "diff",
    "Segmentname_FV2404",
    "Segmentgruppe_FV2404",
    "Segment_FV2404",
    "Datenelement_FV2404",
    "Segment ID_FV2404",
    "Code_FV2404",
    "Qualifier_FV2404",
    "Beschreibung_FV2404",
    "Bedingungsausdruck_FV2404",
    "Bedingung_FV2404"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# --- 2. Create the Excel Table over the used range ---------------------
$rng = $ws.Range("A1:U58")
$tbl = $ws.ListObjects.Add(1, $rng, $null, 1)
$tbl.Name = "Table1"

# --- 3. Freeze the header row -------------------------------------------
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
